# Mark additional TODO-list items as complete (up to 40% of the list),
# add a note on the newly-completed conclusions item, and move the
# visible selection down to the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose COMPLETE (column D) flag flips from FALSE to TRUE.
$rows = @(2, 12, 18, 31, 39, 40, 55)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = $true
}

# Row 40 (Conclusions: Link results back to research questions) also
# gets a note in column E explaining why it is now complete.
$ws.Cells.Item(40, 5).Value = "Already was completed"

# Move the view / active selection further down the sheet.
$ws.Range("D58").Select()
